$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.452.36'
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").Value = '2.243.48'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = "'245.18"
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").Value = "'0.629"
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("D7").Value = "'75.41"
$ws.Range("E7").Value = '  -2.18%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = "'0.621"
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").Value = "'43.57"
$ws.Range("E10").Value = '  +5.87%  '
$ws.Range("D11").Value = "'0.0949"
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("D14").Value = "'14.54"
$ws.Range("E14").Value = '  -2.26%  '
$ws.Range("D15").Value = "'0.858"
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").Value = '2.233.69'
$ws.Range("E16").Value = '  -0.65%  '
$ws.Range("D17").Value = '42.293.10'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("E18").Value = '  +4.12%  '
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("D20").Value = "'72.00"
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = "'10.73"
$ws.Range("E21").Value = '  +47.65%  '
$ws.Range("E22").Value = '  -4.99%  '
$ws.Range("D23").Value = "'231.59"
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = "'11.67"
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  -1.25%  '
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = '  +4.45%  '
$ws.Range("D29").Value = "'167.03"
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").Value = "'20.74"
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("E31").Value = '  +25.01%  '
$ws.Range("D32").Value = "'0.0817"
$ws.Range("E32").Value = '  -1.98%  '
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("D34").Value = "'30.47"
$ws.Range("E34").Value = '  -7.98%  '
$ws.Range("D35").Value = "'0.125"
$ws.Range("E35").Value = '  +0.63%  '
$ws.Range("E36").Value = '  +3.04%  '
$ws.Range("E37").Value = '  +4.36%  '
$ws.Range("D38").Value = "'13.45"
$ws.Range("E38").Value = '  -5.33%  '
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("D40").Value = "'5.71"
$ws.Range("E40").Value = '  -3.60%  '
$ws.Range("D41").Value = "'63.57"
$ws.Range("E41").Value = '  +3.41%  '
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("D43").Value = "'106.59"
$ws.Range("E43").Value = '  -4.77%  '
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("E45").Value = '  +1.78%  '
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = '  +6.81%  '
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("D50").Value = "'4.16"
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("E51").Value = '  +1.30%  '
